# Update loading_percent values for the "case with 380 kV done" commit.
# Only columns B,C,D,E,G,I,J,K,L,M for rows 2-25 change; F,H,N,O,A stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","I","J","K","L","M")

$data = @(
    @(21.63680145232487,12.1619385193875,12.53841595259499,11.41297947709919,3.904430652023717,62.72879921337066,7.243225392641117,20.77577049963612,14.64962145754878,21.41411060295107),
    @(21.75076131927398,11.97208812635933,12.55384038386518,11.43786491540894,3.908580487355095,61.63273878443471,7.236845207378141,20.82669145876357,14.6574557958328,21.48356447772231),
    @(21.82801760390798,11.85595966337181,12.56538285881444,11.45411273770534,3.911255819470507,60.94673472976633,7.232856925650043,20.86367594471317,14.6639616645725,21.53033783491196),
    @(21.86132405356248,11.80879741261374,12.57060723956846,11.46097798719646,3.912378189355412,60.66410854421592,7.23121376520335,20.88018094927612,14.6670394494302,21.55043563516628),
    @(21.86696449167817,11.80097735219987,12.57150618576614,11.46213272218223,3.912566503793927,60.61699897494812,7.230939835576406,20.88300803518332,14.66757628371041,21.55383548362649),
    @(21.82845941207459,11.85532289824388,12.56545120862508,11.45420433545891,3.911270825791076,60.94293531245862,7.232834838027554,20.86389273902239,14.66400144515624,21.5306046821806),
    @(21.67457822030063,12.09641683274433,12.54330416862887,11.42135947739801,3.905835172569901,62.35369733037022,7.2410401459104,20.79213849503649,14.65197097891904,21.43720094952474),
    @(21.43096590277708,12.57049955950812,12.51632157048008,11.36460037047252,3.896179873800516,65.00849353974999,7.256573557652494,20.69699276591504,14.6418218791825,21.28683139964363),
    @(21.28790729342845,12.91670345255705,12.50653313934361,11.32751931715049,3.889689495274411,66.880181086017,7.267653630233573,20.65510211853697,14.64254461553061,21.19639976140359),
    @(21.23072668120876,13.07315032705637,12.50425923004904,11.31164414188449,3.886866004047512,67.71245359176069,7.272623689201124,20.64216983310566,14.64464483363637,21.15962358061183),
    @(21.21021759964756,13.13219507257375,12.50371125436197,11.30577473489677,3.885815229153328,68.02469299511523,7.274495718417127,20.63815601478177,14.6456942952899,21.14632548260727),
    @(21.21458358758039,13.11948829884318,12.503815349993,11.30703250348178,3.886040715197688,67.95757927590003,7.274092987866101,20.63898113624968,14.64545698135291,21.14916150826519),
    @(21.2290164125703,13.07801221802306,12.50420787432732,11.31115841644589,3.886779187758004,67.73820124699591,7.272777902569666,20.64182189606887,14.64472608341736,21.15851694126031),
    @(21.23800615869476,13.05257969509678,12.50448907387152,11.31370415402452,3.88723391842464,67.60344024207296,7.271971072003944,20.64367705584487,14.64431146733515,21.1643292605218),
    @(21.29180381334193,12.90645412296134,12.50672557352182,11.32857672763091,3.889876601866278,66.82538985460508,7.267327417664321,20.65607074861505,14.64244295975583,21.19889102852091),
    @(21.32683563269461,12.81650874500025,12.50865554135057,11.33795448971448,3.891530752186312,66.34304943972137,7.264460737753726,20.66524447180388,14.6417500138921,21.22121142454469),
    @(21.3477278659952,12.76467915244355,12.50997071605211,11.34344185141988,3.892494327053639,66.06382529490503,7.262805275054244,20.67109723868179,14.64151822566452,21.23446000205558),
    @(21.35492894748175,12.74711566071784,12.51045124396249,11.34531585973499,3.892822667980687,65.96898182640911,7.262243623427139,20.67317777113909,14.64146840636268,21.23901620997746),
    @(21.32302949204742,12.82609381219631,12.50842886628847,11.33694653652797,3.891353408400076,66.39458239054147,7.264766585669151,20.66420824880956,14.6418065224245,21.21879288990705),
    @(21.22474602917885,13.09020050580623,12.50408408514499,11.30994268222421,3.886561781607022,67.80271854253898,7.273164446355749,20.64096350250992,14.64493387315208,21.15575196525722),
    @(21.16718382650514,13.26163263300686,12.50306935140294,11.29312252712238,3.883537480833822,68.70590173026285,7.278594545355918,20.63092153347362,14.64845886333791,21.11821342546606),
    @(21.19729268056617,13.17025924063464,12.50344406262174,11.30202416992092,3.885141832361216,68.2254748939912,7.275701712908286,20.63580913138224,14.64644220070778,21.13791300704769),
    @(21.32474790663565,12.8217607696577,12.5085307056376,11.33740193321648,3.891433546325729,66.37129031167935,7.264628334789823,20.66467492300623,14.6417804559659,21.21988501286267),
    @(21.49059497013924,12.44241122149737,12.52185893175084,11.37914079284232,3.898685302382134,64.30339493865202,7.25243024607187,20.71782827981118,14.64312968613135,21.32399394975226)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}
